# Word Replace and Send Email Complete
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Word replace: JR30744 was re-issued as JR73297 ---
$ws.Cells.Replace("JR30744", "JR73297")

# --- Fill in Address / City / State / Zip columns pulled from the new hire records ---

# Row 2 - JR22644 / Annapolis Junction, MD
$ws.Range("F2").Value = "087 Lake Floyd Circle"
$ws.Range("G2").Value = "Annapolis Junction"
$ws.Range("H2").Value = "Maryland"
$ws.Range("I2").Value = 20701

# Row 3 - JR22755 / Bonita Springs, FL
$ws.Range("F3").Value = "185 Wilkinson Court"
$ws.Range("G3").Value = "Bonita Springs"
$ws.Range("H3").Value = "Florida"
$ws.Range("I3").Value = 33293

# Row 4 - JR73297 / Bothell, WA
$ws.Range("F4").Value = "219 Ryder Avenue"
$ws.Range("G4").Value = "Bothell"
$ws.Range("H4").Value = "Washinton"
$ws.Range("I4").Value = 98011

# Row 5 - JR12356 / West Chester, PA
$ws.Range("F5").Value = "860 Davis Avenue"
$ws.Range("G5").Value = "West Chester"
$ws.Range("H5").Value = "Pennsylvania"
$ws.Range("I5").Value = 19382

# Row 6 - JR73297 / Tampa, FL
$ws.Range("F6").Value = "813 Collins Street"
$ws.Range("G6").Value = "Tampa"
$ws.Range("H6").Value = "Florida"
$ws.Range("I6").Value = 33634

# Row 7 - JR73297 / Borentown, NJ
$ws.Range("F7").Value = "2267 Lake Road"
$ws.Range("G7").Value = "Borentown"
$ws.Range("H7").Value = "New Jersey"
$ws.Range("I7").Value = 85050

# --- Formatting: the City cell pasted for row 3 came in with the web source's font ---
$ws.Range("G3").Font.Name = "Open Sans"
$ws.Range("G3").Font.Size = 10.5
$ws.Range("G3").Font.Color = 3355443

# Row grew slightly taller to fit the new font
$ws.Rows.Item(3).RowHeight = 15.75

# --- Move the active selection to reflect where the user left off (just past the table) ---
$ws.Range("I8").Select()

Write-Output "Edit complete"
